$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the remark in F3 ("actually venous pH for now") ---
$ws.Range("F3").ClearContents()

# --- Unmerge row 29 (A29:E29) and fill it in with the new "thrombocytes" row ---
$ws.Range("A29:E29").UnMerge()

$ws.Range("A29").Value = "thrombocytes"
$ws.Range("B29").Value = "thrombocytes"
$ws.Range("C29").Value = "G/l"
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 1000

# Match formatting used by the analogous "weight / kg" row (row 27)
$ws.Range("A29").NumberFormat = "General"
$ws.Range("B29").NumberFormat = "General"
$ws.Range("C29").NumberFormat = "0"
$ws.Range("C29").HorizontalAlignment = -4108
$ws.Range("D29").NumberFormat = "0"
$ws.Range("E29").NumberFormat = "#,##0"

# --- Update the active selection to match the new cursor position ---
$ws.Range("B32").Select()
